$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill previously-blank "backup" (column R) values for rows 312-315 with 0 ---
$ws.Range("R312").Value = 0
$ws.Range("R313").Value = 0
$ws.Range("R314").Value = 0
$ws.Range("R315").Value = 0
$ws.Range("R316").Value = 0

# --- 2. Row 314 "isPivot" (column O) flips from 0 to 1 ---
$ws.Range("O314").Value = 1

# --- 3. Append new daily rows 317-324 (stock.yaml completed backfill) ---
$newRows = @(
    @(317, 45635, 338.5,              341.2999877929688, 335.8999938964844, 336.8500061035156, 336.8500061035156, 10589890, 2024, 12, 9,  0, 0, 0, 50, 0, 0, 0),
    @(318, 45636, 337,                337.9500122070312, 333.6499938964844, 335.1000061035156, 335.1000061035156, 7366133,  2024, 12, 10, 0, 0, 0, 50, 0, 0, 0),
    @(319, 45637, 334.7000122070312,  345,                333.2999877929688, 343,                343,                20098920, 2024, 12, 11, 0, 0, 0, 50, 0, 0, 0),
    @(320, 45638, 343.7999877929688,  346.2999877929688, 337.2000122070312, 338.2999877929688, 338.2999877929688, 10639640, 2024, 12, 12, 0, 0, 0, 50, 0, 0, 0),
    @(321, 45639, 337,                340.7000122070312, 330.6499938964844, 339.75,             339.75,             9126606,  2024, 12, 13, 0, 0, 0, 50, 0, 0, 0),
    @(322, 45642, 339.25,             342.4500122070312, 338.1499938964844, 340.3500061035156, 340.3500061035156, 6632276,  2024, 12, 16, 0, 0, 0, 51, 0, 0, 0),
    @(323, 45643, 339.0499877929688,  340.8999938964844, 335.8999938964844, 337.1000061035156, 337.1000061035156, 8027327,  2024, 12, 17, 0, 0, 0, 51, 0, 0, 0),
    @(324, 45644, 336,                337,                323.5499877929688, 324.6000061035156, 324.6000061035156, 13945890, 2024, 12, 18, 0, 0, 0, 51, 0, 0, 0)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # Column A (Datetime) keeps the same date/time number format as the rows above it
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # Columns B..Q (Open .. detect_structure)
    for ($c = 2; $c -le 17; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
    # Column R ("backup") is intentionally left blank for the new rows,
    # matching the source data (still unset / inline-string empty).
}
